$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.920.97"
$ws.Range("E2").Value = "  +0.52%  "
$ws.Range("D3").Value = "1.632.32"
$ws.Range("E3").Value = "  +1.82%  "
$ws.Range("E4").Value = "  +0.31%  "
$ws.Range("D5").Value = "'214.38"
$ws.Range("E5").Value = "  +0.85%  "
$ws.Range("E6").Value = "  -0.07%  "
$ws.Range("E7").Value = "  +0.25%  "
$ws.Range("D8").Value = "'28.46"
$ws.Range("E8").Value = "  -0.27%  "
$ws.Range("E9").Value = "  +1.11%  "
$ws.Range("D10").Value = "'0.0607"
$ws.Range("E10").Value = "  +0.58%  "
$ws.Range("D11").Value = "'0.0907"
$ws.Range("E11").Value = "  +0.06%  "
$ws.Range("D12").Value = "1.867.28"
$ws.Range("E12").Value = "  +1.91%  "
$ws.Range("D13").Value = "1.645.86"
$ws.Range("E13").Value = "  +2.96%  "
$ws.Range("D14").Value = "'0.562"
$ws.Range("E14").Value = "  +1.92%  "
$ws.Range("D15").Value = "'9.24"
$ws.Range("E15").Value = "  +15.57%  "
$ws.Range("D16").Value = "29.939.80"
$ws.Range("E16").Value = "  +0.71%  "
$ws.Range("D17").Value = "'3.84"
$ws.Range("E17").Value = "  +1.75%  "
$ws.Range("D18").Value = "'63.98"
$ws.Range("E18").Value = "  -0.08%  "
$ws.Range("D19").Value = "'241.60"
$ws.Range("E19").Value = "  -0.32%  "
$ws.Range("D20").Value = "0.0₃0700"
$ws.Range("E20").Value = "  +0.22%  "
$ws.Range("E21").Value = "  +0.16%  "
$ws.Range("D22").Value = "'4.13"
$ws.Range("E22").Value = "  +2.11%  "
$ws.Range("D23").Value = "'9.76"
$ws.Range("E23").Value = "  +3.16%  "
$ws.Range("E24").Value = "  +2.91%  "
$ws.Range("D25").Value = "'158.40"
$ws.Range("E25").Value = "  +2.02%  "
$ws.Range("D26").Value = "'15.50"
$ws.Range("E26").Value = "  +0.21%  "
$ws.Range("E27").Value = "  +0.33%  "
$ws.Range("D28").Value = "'6.60"
$ws.Range("E28").Value = "  +2.28%  "
$ws.Range("E29").Value = "  +0.30%  "
$ws.Range("E30").Value = "  +1.81%  "
$ws.Range("E31").Value = "  +4.03%  "
$ws.Range("E32").Value = "  +3.91%  "
$ws.Range("D33").Value = "'3.17"
$ws.Range("E33").Value = "  -0.48%  "
$ws.Range("D34").Value = "1.423.88"
$ws.Range("E34").Value = "  +0.08%  "
$ws.Range("D35").Value = "'1.64"
$ws.Range("E35").Value = "  +4.51%  "
$ws.Range("E36").Value = "  -0.91%  "
$ws.Range("D37").Value = "'2.79"
$ws.Range("E37").Value = "  -2.58%  "
$ws.Range("D38").Value = "'2.30"
$ws.Range("E38").Value = "  -0.23%  "
$ws.Range("E39").Value = "  -0.15%  "
$ws.Range("D40").Value = "'75.52"
$ws.Range("E40").Value = "  +12.41%  "
$ws.Range("E41").Value = "  +1.46%  "
$ws.Range("E42").Value = "  +3.16%  "
$ws.Range("D43").Value = "'0.827"
$ws.Range("E43").Value = "  +1.12%  "
$ws.Range("E45").Value = "  +2.25%  "
$ws.Range("E46").Value = "  +0.30%  "
$ws.Range("D47").Value = "'52.78"
$ws.Range("E47").Value = "  -4.29%  "
$ws.Range("D48").Value = "'5.35"
$ws.Range("E48").Value = "  -0.60%  "
$ws.Range("D49").Value = "1.774.14"
$ws.Range("E49").Value = "  +1.95%  "
$ws.Range("E50").Value = "  +9.14%  "
$ws.Range("D51").Value = "'90.50"
$ws.Range("E51").Value = "  +4.49%  "
